# This script inserts a new row of data at row 2 of the sheet ("南京希音电子商务有限公司"),
# pushing the previously-existing rows 2-12 down to rows 3-13, while keeping the
# column A sequence numbers (0,1,2,...) fixed per row position (row 2 stays 0, row 3 stays 1, etc.)
# and appending a new trailing index (11) for the row that lands at row 13.
#
# Because Range.Value cannot be reliably read back in this runtime, the row content is
# shifted down using Range.Copy (which copies actual cell contents/formatting), after first
# clearing the destination so that cells that should become blank do not retain stale data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: create new row 13 from current row 12's content (columns B:P) ---
$ws.Range("B13:P13").ClearContents()
$ws.Range("B12:P12").Copy($ws.Range("B13:P13"))

# Column A13: copy style from A12, then set its own sequence value (11)
$ws.Range("A12").Copy($ws.Range("A13"))
$ws.Range("A13").Value = 11

# --- Step 2: shift rows 11..2 down into rows 12..3 (columns B:P) ---
for ($r = 12; $r -ge 3; $r--) {
    $src = $r - 1
    $ws.Range("B$r`:P$r").ClearContents()
    $ws.Range("B$src`:P$src").Copy($ws.Range("B$r`:P$r"))
}

# --- Step 3: write the brand-new row 2 content ---
$ws.Range("B2").Value = "南京希音电子商务有限公司"
$ws.Range("C2").Value = "天溯产业园"
$ws.Range("D2").Value = $null
$ws.Range("E2").Value = "前端"
$ws.Range("F2").Value = "10:00-18:00(到20:00有50补贴)"
$ws.Range("G2").Value = "12:00-13:30"
$ws.Range("H2").Value = "看部门，不强制，有工时排名。"
$ws.Range("I2").Value = "基础工资的8%"
$ws.Range("J2").Value = "看部门盈利情况和个人绩效定"
$ws.Range("K2").Value = "试用期6个月，100%工资不打折"
$ws.Range("L2").Value = "配mac m1+显示器，网吧工作环境，工位挤。"
$ws.Range("M2").Value = "法定年假，可用加班时长来调休"
$ws.Range("N2").Value = "1月3次补卡"
$ws.Range("O2").Value = "抠，舍得给校招生，不舍得给社招生。多余的调休时长换钱200/d"
$ws.Range("P2").Value = $null
